$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# Row 2 (married)
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 61.11111111111111
$ws.Range("D2").Value = 33
$ws.Range("E2").Value = 44.5945945945946

# Row 3 (single)
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 38.88888888888889
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 50

# Row 4 - values moved from B4/C4 to D4/E4
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 5.405405405405405
